$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2028216.2
$ws.Range("I137").Value = 5918169
$ws.Range("J137").Value = 5440.84
$ws.Range("K137").Value = 17754507
$ws.Range("L137").Value = 16322.52
$ws.Range("M137").Value = -17751957
$ws.Range("N137").Value = -21422.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2177.077
$ws.Range("I61").Value = 1431.3158
$ws.Range("J61").Value = 4201.2856
$ws.Range("K61").Value = 1431.3158
$ws.Range("L61").Value = 4201.2856
$ws.Range("M61").Value = -1219.3158
$ws.Range("N61").Value = -4625.2856
$ws.Range("H74").Value = 2003.683
$ws.Range("I74").Value = 1718.4333
$ws.Range("J74").Value = 2781.6365
$ws.Range("K74").Value = 1718.4333
$ws.Range("L74").Value = 2781.6365
$ws.Range("M74").Value = -844.4332999999999
$ws.Range("N74").Value = -4529.636500000001
$ws.Range("H77").Value = 2003.683
$ws.Range("I77").Value = 1718.4333
$ws.Range("J77").Value = 2781.6365
$ws.Range("K77").Value = 8592.166499999999
$ws.Range("L77").Value = 13908.1825
$ws.Range("M77").Value = -4224.166499999999
$ws.Range("N77").Value = -22644.1825
$ws.Range("H132").Value = 16668555
$ws.Range("I132").Value = 20834486
$ws.Range("J132").Value = 4832.3335
$ws.Range("K132").Value = 62503458
$ws.Range("L132").Value = 14497.0005
$ws.Range("M132").Value = -62500928
$ws.Range("N132").Value = -19557.0005
$ws.Range("H136").Value = 2177.077
$ws.Range("I136").Value = 1431.3158
$ws.Range("J136").Value = 4201.2856
$ws.Range("K136").Value = 4293.9474
$ws.Range("L136").Value = 12603.8568
$ws.Range("M136").Value = -1743.9474
$ws.Range("N136").Value = -17703.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1263.2424
$ws.Range("I94").Value = 1178.2084
$ws.Range("J94").Value = 1490
$ws.Range("K94").Value = 1178.2084
$ws.Range("L94").Value = 1490
$ws.Range("M94").Value = -727.2084
$ws.Range("N94").Value = -2392
$ws.Range("H134").Value = 2803.327
$ws.Range("I134").Value = 1606.1
$ws.Range("J134").Value = 3551.5938
$ws.Range("K134").Value = 4818.299999999999
$ws.Range("L134").Value = 10654.7814
$ws.Range("M134").Value = -2283.299999999999
$ws.Range("N134").Value = -15724.7814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5169.468
$ws.Range("I31").Value = 2222.6924
$ws.Range("J31").Value = 5951.265
$ws.Range("K31").Value = 2222.6924
$ws.Range("L31").Value = 5951.265
$ws.Range("M31").Value = -1927.6924
$ws.Range("N31").Value = -6541.265
$ws.Range("H34").Value = 5169.468
$ws.Range("I34").Value = 2222.6924
$ws.Range("J34").Value = 5951.265
$ws.Range("K34").Value = 2222.6924
$ws.Range("L34").Value = 5951.265
$ws.Range("M34").Value = -2020.6924
$ws.Range("N34").Value = -6355.265
$ws.Range("H58").Value = 1948.7354
$ws.Range("I58").Value = 1376.8948
$ws.Range("J58").Value = 2673.0667
$ws.Range("K58").Value = 1376.8948
$ws.Range("L58").Value = 2673.0667
$ws.Range("M58").Value = -1173.8948
$ws.Range("N58").Value = -3079.0667
$ws.Range("H132").Value = 52150.43
$ws.Range("I132").Value = 1637.4706
$ws.Range("J132").Value = 130215.91
$ws.Range("K132").Value = 4912.4118
$ws.Range("L132").Value = 390647.73
$ws.Range("M132").Value = -2382.4118
$ws.Range("N132").Value = -395707.73
$ws.Range("H134").Value = 370173.97
$ws.Range("I134").Value = 381978.7
$ws.Range("K134").Value = 1145936.1
$ws.Range("M134").Value = -1143401.1
$ws.Range("H136").Value = 1948.7354
$ws.Range("I136").Value = 1376.8948
$ws.Range("J136").Value = 2673.0667
$ws.Range("K136").Value = 4130.6844
$ws.Range("L136").Value = 8019.2001
$ws.Range("M136").Value = -1580.6844
$ws.Range("N136").Value = -13119.2001
$ws.Range("H140").Value = 31994.5
$ws.Range("J140").Value = 31994.5
$ws.Range("L140").Value = 31994.5
$ws.Range("N140").Value = -42354.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 927.9
$ws.Range("I98").Value = 450
$ws.Range("J98").Value = 981
$ws.Range("K98").Value = 1350
$ws.Range("L98").Value = 2943
$ws.Range("M98").Value = 148
$ws.Range("N98").Value = -5939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20836802
$ws.Range("I132").Value = 32261092
$ws.Range("J132").Value = 4276.5293
$ws.Range("K132").Value = 96783276
$ws.Range("L132").Value = 12829.5879
$ws.Range("M132").Value = -96780746
$ws.Range("N132").Value = -17889.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2214.4517
$ws.Range("I7").Value = 1637.2084
$ws.Range("J7").Value = 4193.5713
$ws.Range("K7").Value = 1637.2084
$ws.Range("L7").Value = 4193.5713
$ws.Range("M7").Value = -1525.2084
$ws.Range("N7").Value = -4417.5713
$ws.Range("H93").Value = 1683.6666
$ws.Range("I93").Value = 800
$ws.Range("K93").Value = 800
$ws.Range("M93").Value = 448
$ws.Range("H126").Value = 2214.4517
$ws.Range("I126").Value = 1637.2084
$ws.Range("J126").Value = 4193.5713
$ws.Range("K126").Value = 4911.6252
$ws.Range("L126").Value = 12580.7139
$ws.Range("M126").Value = -2441.6252
$ws.Range("N126").Value = -17520.7139
$ws.Range("H132").Value = 2837.879
$ws.Range("I132").Value = 2186.923
$ws.Range("J132").Value = 5255.7144
$ws.Range("K132").Value = 6560.768999999999
$ws.Range("L132").Value = 15767.1432
$ws.Range("M132").Value = -4030.768999999999
$ws.Range("N132").Value = -20827.1432
$ws.Range("H136").Value = 1666.75
$ws.Range("I136").Value = 1216.9
$ws.Range("J136").Value = 3916
$ws.Range("K136").Value = 3650.7
$ws.Range("L136").Value = 11748
$ws.Range("M136").Value = -1100.7
$ws.Range("N136").Value = -16848

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1176594
$ws.Range("I132").Value = 1403551.8
$ws.Range("J132").Value = 3979
$ws.Range("K132").Value = 4210655.4
$ws.Range("L132").Value = 11937
$ws.Range("M132").Value = -4208125.4
$ws.Range("N132").Value = -16997
$ws.Range("H136").Value = 467742.6
$ws.Range("I136").Value = 707982.6
$ws.Range("J136").Value = 1394.2354
$ws.Range("K136").Value = 2123947.8
$ws.Range("L136").Value = 4182.706200000001
$ws.Range("M136").Value = -2121397.8
$ws.Range("N136").Value = -9282.706200000001
